$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at C and D for Start Date / End Date (shifts old C:F -> E:H)
$ws.Columns("C:D").Insert()
$ws.Columns("C:D").ColumnWidth = 21.1796875

# Rename the payment terms header (now shifted into column F after the insert)
$ws.Range("F1").Value = "Payment_Terms_Contract"

# New header labels for the inserted columns
$ws.Range("C1").Value = "Start Date"
$ws.Range("D1").Value = "End Date"

# Fill in the date values (as raw serial numbers so Excel doesn't pre-apply
# its own guessed date format) for the two data rows
$ws.Range("C2").Value = 41849
$ws.Range("D2").Value = 44196
$ws.Range("C3").Value = 41849
$ws.Range("D3").Value = 44196

# Apply short-date number format to C2, then copy that formatting across
# C2:D5 so every date cell shares the same style (matches native Excel
# fill/format-paint behaviour instead of creating a distinct style per cell).
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy()
$ws.Range("C2:D5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
